$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text format
# so Excel stores them as the literal string (matching the source data) instead of
# silently converting them into numeric values and dropping formatting like trailing zeros.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Apply the updated coin price / link / volume figures
$ws.Range("D2").Value = "30.143.89"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.830.31"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "231.06"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4651"
$ws.Range("E7").Value = "  -3.72%  "
$ws.Range("D8").Value = "0.2685"
$ws.Range("E8").Value = "  -6.65%  "
$ws.Range("D9").Value = "0.06260"
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("D10").Value = "1.852.13"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "0.07384"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "15.99"
$ws.Range("E12").Value = "  -4.69%  "
$ws.Range("D13").Value = "4.895"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").Value = "83.11"
$ws.Range("E14").Value = "  -5.41%  "
$ws.Range("D15").Value = "0.6173"
$ws.Range("E15").Value = "  -7.34%  "
$ws.Range("D16").Value = "30.065.92"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "226.27"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "0.000007269"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").Value = "12.36"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "2.071.70"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "4.833"
$ws.Range("E23").Value = "  -8.37%  "
$ws.Range("B24").Value = "BitDAO"
$ws.Range("C24").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D24").Value = "0.3911"
$ws.Range("E24").Value = "  +10.32%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "5.852"
$ws.Range("E25").Value = "  -5.42%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.109"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "164.38"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "17.59"
$ws.Range("E28").Value = "  -6.00%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.842"
$ws.Range("E29").Value = "  -5.85%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.1007"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "1.367"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.755"
$ws.Range("E33").Value = "  -6.67%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04769"
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6996"
$ws.Range("E36").Value = "  -6.82%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01808"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.610"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.8909"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "1.921"
$ws.Range("E41").Value = "  -6.94%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "102.62"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.462"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3990"
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.939"
$ws.Range("E46").Value = "  -6.47%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1188"
$ws.Range("E47").Value = "  -6.98%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "59.48"
$ws.Range("E48").Value = "  -7.32%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.448"
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("D50").Value = "32.57"
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05519"
$ws.Range("E51").Value = "  -2.63%  "
